# Apply "Calculations for Linear Voltage Regulators" to the eagle
# (Spannungsregler) sheet: add a small "R1 unten" / "R2 oben" label pair
# in columns G1:G2, and make the Spannungsregler sheet the active tab.

$wb = $excel.ActiveWorkbook

$wsReg = $wb.Worksheets.Item("Spannungsregler")
$wsBatt = $wb.Worksheets.Item("Batteriespannungsmessung")

# Add the two new labels used when eagle-computing the divider resistors.
$wsReg.Range("G1").Value = "R1 unten"
$wsReg.Range("G2").Value = "R2 oben"

# Select G3 on the Spannungsregler sheet (matches the selection recorded
# in the workbook after the edit) and activate that sheet/cell.
$wsReg.Activate()
$wsReg.Range("G3").Select()

# The Batteriespannungsmessung sheet keeps its own remembered selection,
# but it is no longer the tab that is selected/active when the file is
# reopened.
$wsBatt.Range("B6").Select()

# Re-activate Spannungsregler last so it becomes the workbook's active
# sheet/tab.
$wsReg.Activate()
